# V1.1 Tighter Tolerance for smaller CF tubes
#
# The "HD_CFx_MGN9_CFTube Insert" rows (19-22, column J "Comment") get a new
# note explaining that different CF-tube insert files should be used
# depending on the measured tube tolerance.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$comment = "For CF tube within 20.3-20.5mm use V1 Files" + [char]10 + "For CF Tubes within 20.1-20.2mm use V1.1 Files"

$rows = 19, 20, 21, 22
foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 10)   # column J = Comment
    $cell.Value = $comment
    $cell.WrapText = $true
    $cell.HorizontalAlignment = -4131   # xlLeft
    $cell.VerticalAlignment = -4108     # xlCenter
}

# Reflect the author's final cursor position/selection in the saved view.
$ws.Range("I19").Select() | Out-Null
